# Apply updated crypto market data (prices and 1h volume % changes) to the sheet.
# Matches the commit "Updated cryptos list on Sat Apr 13 04:45:00 UTC 2024 with GitHub Actions".
#
# Cell values in columns D (Price) and E (Volume(1h)) are stored as plain text
# (not numbers) in the source data -- some price strings (e.g. "67.429.67",
# "0.998", "26.24") look numeric to Excel's auto-detection, so a leading
# apostrophe is used to force text entry, then the cell style is reset back
# to "Normal" so no stray quote-prefix / text-number-format style lingers on
# the cell (keeping it identical in shape to the original workbook's cells,
# which carry no explicit style).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cellRef, $value) {
    $ws.Range($cellRef).Value = "'" + $value
    $ws.Range($cellRef).Style = "Normal"
}

Set-TextValue "D2" "67.429.67"
Set-TextValue "E2" "  -4.95%  "
Set-TextValue "D3" "3.238.57"
Set-TextValue "E3" "  -8.58%  "
Set-TextValue "D4" "0.998"
Set-TextValue "E4" "  -0.12%  "
Set-TextValue "D5" "585.79"
Set-TextValue "E5" "  -5.31%  "
Set-TextValue "D6" "152.99"
Set-TextValue "E6" "  -12.00%  "
Set-TextValue "D7" "0.999"
Set-TextValue "E7" "  -0.15%  "
Set-TextValue "D8" "3.233.75"
Set-TextValue "E8" "  -8.60%  "
Set-TextValue "E9" "  -10.94%  "
Set-TextValue "E10" "  -12.48%  "
Set-TextValue "D11" "6.77"
Set-TextValue "E11" "  -5.91%  "
Set-TextValue "D12" "0.506"
Set-TextValue "E12" "  -14.09%  "
Set-TextValue "D13" "38.54"
Set-TextValue "E13" "  -17.52%  "
Set-TextValue "D14" "0.0000245"
Set-TextValue "E14" "  -11.42%  "
Set-TextValue "D15" "3.757.91"
Set-TextValue "E15" "  -8.65%  "
Set-TextValue "D16" "67.328.48"
Set-TextValue "E16" "  -5.10%  "
Set-TextValue "D17" "3.236.95"
Set-TextValue "E17" "  -8.33%  "
Set-TextValue "D18" "543.49"
Set-TextValue "E18" "  -11.32%  "
Set-TextValue "E19" "  -5.78%  "
Set-TextValue "D20" "7.21"
Set-TextValue "E20" "  -14.82%  "
Set-TextValue "D21" "15.22"
Set-TextValue "E21" "  -14.64%  "
Set-TextValue "D22" "0.764"
Set-TextValue "E22" "  -14.15%  "
Set-TextValue "D23" "7.81"
Set-TextValue "E23" "  -13.60%  "
Set-TextValue "D24" "85.92"
Set-TextValue "E24" "  -12.60%  "
Set-TextValue "D25" "13.56"
Set-TextValue "E25" "  -13.89%  "
Set-TextValue "E26" "  +0.06%  "
Set-TextValue "D27" "3.18"
Set-TextValue "E27" "  -16.28%  "
Set-TextValue "D28" "8.15"
Set-TextValue "E28" "  -11.04%  "
Set-TextValue "D29" "29.52"
Set-TextValue "E29" "  -12.84%  "
Set-TextValue "D30" "2.15"
Set-TextValue "E30" "  -17.53%  "
Set-TextValue "D31" "2.70"
Set-TextValue "E31" "  -11.34%  "
Set-TextValue "D32" "1.16"
Set-TextValue "E32" "  -11.24%  "
Set-TextValue "D33" "543.84"
Set-TextValue "E33" "  -10.23%  "
Set-TextValue "D34" "6.60"
Set-TextValue "E34" "  -19.49%  "
Set-TextValue "D35" "5.77"
Set-TextValue "E35" "  -16.11%  "
Set-TextValue "E36" "  +0.06%  "
Set-TextValue "D37" "0.0450"
Set-TextValue "E37" "  -5.71%  "
Set-TextValue "D38" "53.85"
Set-TextValue "E38" "  -5.50%  "
Set-TextValue "D39" "0.0852"
Set-TextValue "E39" "  -15.40%  "
Set-TextValue "D40" "9.22"
Set-TextValue "E40" "  -15.06%  "
Set-TextValue "E41" "  -12.96%  "
Set-TextValue "D42" "2.932.72"
Set-TextValue "E42" "  -13.29%  "
Set-TextValue "D43" "2.64"
Set-TextValue "E43" "  -25.16%  "
Set-TextValue "D44" "0.0₃0591"
Set-TextValue "E44" "  -20.26%  "
Set-TextValue "D45" "0.262"
Set-TextValue "E45" "  -16.77%  "
Set-TextValue "D46" "2.40"
Set-TextValue "E46" "  -20.19%  "
Set-TextValue "B47" "InjectiveProtocol"
Set-TextValue "C47" "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
Set-TextValue "D47" "26.24"
Set-TextValue "E47" "  -18.84%  "
Set-TextValue "B48" "USDe"
Set-TextValue "C48" "https://coinranking.com/coin/exbfr2U-0+usde-usde"
Set-TextValue "D48" "1.00"
Set-TextValue "E48" "  -0.05%  "
Set-TextValue "D49" "2.13"
Set-TextValue "E49" "  -17.13%  "
Set-TextValue "D50" "0.114"
Set-TextValue "E50" "  -13.15%  "
Set-TextValue "D51" "124.58"
Set-TextValue "E51" "  -6.74%  "
